$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameters")

# Update the shared string text used for the "Hoffman et al., 2013" citation
# to also credit Clegg et al., 2017 (affects E27 and E28, which both cite it)
$ws.Range("E27").Value = "Clegg et al., 2017 & Hoffman et al., 2013"
$ws.Range("E28").Value = "Clegg et al., 2017 & Hoffman et al., 2013"

# Fill in the previously-empty association/dissociation rate values for
# PlGF:NRP1 binding
$ws.Range("C27").Value = 10000
$ws.Range("C28").Value = 0.001

# Update the selected cell in the sheet view to C28
$ws.Range("C28").Select()
